$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rule statements (column B) and its description (column A) for the
# "check whether a password changed in the last 7 days" scenario.
$statements = @"
statements:
  - context: pwdChanged := user["urn:ietf:params:scim:schemas:extension:ibm:2.0:User"].pwdChangedTime
  - context: currentTime := now
  - context: timeDiff := timestamp(context.currentTime) - timestamp(context.pwdChanged)
  - if:
        match: context.timeDiff <= duration('604800s')
        block:
            - return: true
  - return: false
"@

$description = "check whether a password changed in the last 7 days"

$ws.Range("B71").Value = $statements
$ws.Range("A71").Value = $description

# Match the look of the other data rows: wrapped text and an auto-fit-style
# row height.
$ws.Range("A71:B71").WrapText = $true
$ws.Rows.Item(71).RowHeight = 129.6

# Update the view so the new row is visible / selected, mirroring the
# author's scroll position and selection after the edit.
$ws.Range("A77").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
